$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.816149353981018
$ws.Range("B1").Value = 4.551989555358887
$ws.Range("C1").Value = 3.946398019790649
$ws.Range("D1").Value = 0.9041378498077393
$ws.Range("E1").Value = 0.4745079874992371
